# Debugging the sensitivity analysis:
# - Append a second cohort tag "2022_06" (with trailing spaces) into the
#   param_sa_cohorts column (R) and add trailing spaces to param_sa_yob (T)
#   for the existing rows (2 and 3).
# - Append two new data rows (4 and 5) for the "2022_06" cohort, mirroring
#   rows 2 and 3 but with updated KCOR/CI values and slope parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while avoiding Excel's automatic
# "looks like a number/date" coercion, and without leaving a lingering
# explicit number format on the cell (it is reset back to Normal/General
# immediately after the value is written).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Update existing rows 2 and 3: R gets ",2022_06   " appended, T gets
# "   " appended (trailing spaces, stored as text).
# ---------------------------------------------------------------------
$ws.Range("R2").Value = "2021_24,2022_06   "
Set-TextValue $ws.Range("T2") "0   "

$ws.Range("R3").Value = "2021_24,2022_06   "
Set-TextValue $ws.Range("T3") "0   "

# ---------------------------------------------------------------------
# Add new row 4 (Dose_num = 1) for cohort 2022_06
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "2022_06"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 44921
$ws.Range("E4").NumberFormat = $ws.Range("E2").NumberFormat()
$ws.Range("F4").Value = 1.118634762881549
$ws.Range("G4").Value = 1.071698934433147
$ws.Range("H4").Value = 1.167626179817872
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 1
Set-TextValue $ws.Range("N4") "4/1/24"
Set-TextValue $ws.Range("O4") "2024-04-01"
$ws.Range("P4").Value = 19
$ws.Range("Q4").Value = 92
$ws.Range("R4").Value = "2021_24,2022_06   "
$ws.Range("S4").Value = "1,0;2,0"
Set-TextValue $ws.Range("T4") "0   "

# ---------------------------------------------------------------------
# Add new row 5 (Dose_num = 2) for cohort 2022_06
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "2022_06"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 44921
$ws.Range("E5").NumberFormat = $ws.Range("E2").NumberFormat()
$ws.Range("F5").Value = 1.060430986706911
$ws.Range("G5").Value = 1.037039328387625
$ws.Range("H5").Value = 1.084350271765076
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 1
Set-TextValue $ws.Range("N5") "4/1/24"
Set-TextValue $ws.Range("O5") "2024-04-01"
$ws.Range("P5").Value = 19
$ws.Range("Q5").Value = 92
$ws.Range("R5").Value = "2021_24,2022_06   "
$ws.Range("S5").Value = "1,0;2,0"
Set-TextValue $ws.Range("T5") "0   "
